$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.6
$ws.Range("G2").Value = 1
$ws.Range("U2").Value = 1
$ws.Range("V2").Value = 1
$ws.Range("W2").Value = 0.5

# Row 3
$ws.Range("C3").Value = 1
$ws.Range("G3").Value = 0.3333333333333333
$ws.Range("U3").Value = 0.6666666666666666
$ws.Range("V3").Value = 0.6666666666666666
$ws.Range("W3").Value = 0.3333333333333333

# Row 4
$ws.Range("C4").Value = 0.7499999999999999
$ws.Range("G4").Value = 0.5
$ws.Range("U4").Value = 0.8
$ws.Range("V4").Value = 0.8
$ws.Range("W4").Value = 0.4

# Row 5
$ws.Range("C5").Value = 0.8823529411764706
$ws.Range("G5").Value = 0.3846153846153846
$ws.Range("U5").Value = 0.7142857142857142
$ws.Range("V5").Value = 0.7142857142857142
$ws.Range("W5").Value = 0.3571428571428571

# Row 6
$ws.Range("C6").Value = 0.9879538239787089
$ws.Range("G6").Value = 0.1064646477465997
$ws.Range("U6").Value = 0.8124242481930318
$ws.Range("V6").Value = 0.5766666455144387
$ws.Range("W6").Value = 0.1064646477465997
